$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The post "「予想外な結末！」" (row 677) was removed from the data set.
# Deleting the entire row shifts all subsequent rows (678-803) up by one,
# matching the shrink of the used range from A1:C803 to A1:C802.
$ws.Rows.Item(677).Delete()
